# Add data for 2022-08-21 (diff shows August counts refreshed through 08-13)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab / internal sheet name
$ws.Name = "Through 2022-08-13"

# Update the "August (through 08-12)" label to "August (through 08-13)"
$ws.Range("A9").Value = "August (through 08-13)"

# Update August row (row 9) values
$ws.Range("B9").Value = 13
$ws.Range("C9").Value = 28
$ws.Range("D9").Value = 28
$ws.Range("G9").Value = 84
$ws.Range("H9").Value = 78
$ws.Range("I9").Value = 75

# Update Total row (row 10) values
$ws.Range("B10").Value = 175
$ws.Range("C10").Value = 330
$ws.Range("D10").Value = 493
$ws.Range("G10").Value = 705
$ws.Range("H10").Value = 988
$ws.Range("I10").Value = 1045
